# "overall offline analysis, maxwellian init option"
#
# - Cameras!C6 boresightAzDeg value tweak: 88.3 -> 88.0172525718237
# - Rows 56 (Sim) / 15 (Cameras) grow a hair taller: 12.8 -> 12.85
# - Selection on Cameras moves from C24 to B15
# - "Arcs" sheet renamed to "Arc0"
# - Arc0 becomes the active/selected sheet (was Sim), with its selection at A1
# - Tab ratio widens slightly (984 -> 988)

$wb = $excel.ActiveWorkbook

$simSheet = $wb.Worksheets.Item("Sim")
$camSheet = $wb.Worksheets.Item("Cameras")
$arcSheet = $wb.Worksheets.Item("Arcs")

# Rename "Arcs" -> "Arc0"
$arcSheet.Name = "Arc0"

# Cameras: new boresightAzDeg value in C6
$camSheet.Range("C6").Value = 88.0172525718237

# Row height bumps
$simSheet.Rows.Item(56).RowHeight = 12.85
$camSheet.Rows.Item(15).RowHeight = 12.85

# Per-sheet selections (applied while each sheet is the active one so each
# sheet keeps its own remembered selection), ending on Arc0 so it is the
# sheet left active/selected for the workbook as a whole.
$simSheet.Activate()
$simSheet.Range("A56").Select()

$camSheet.Activate()
$camSheet.Range("B15").Select()

# Arc0's own (bottom-left / frozen-pane) selection stays put at L4 - only
# its header pane's scroll position moves (W1 -> A1), which activating it
# here does not disturb.
$arcSheet.Activate()

# Widen the sheet-tabs/horizontal-scrollbar split slightly and make sure
# Arc0 (index 2, 0-based) is recorded as the active tab.
$win = $excel.ActiveWindow
$win.TabRatio = 0.988
